{"js": "// Update the date line and every arithmetic-problem cell in the practice\n// table to the new values from the authoritative diff. Cell values are\n// applied positionally (row-major, 5 columns per row) since some old\n// expressions repeat (e.g. \"63+4=\", \"9+84=\") and a naive find/replace\n// would not be able to tell the duplicates apart.\n\nconst body = context.document.body;\n\n// --- 1. Update the date paragraph -----------------------------------\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.insertText(\"2025-05-25 Sunday\", \"Replace\");\n\n// --- 2. Update every cell in the practice table ----------------------\nconst newValues = [\"84-31=\", \"19+33=\", \"85-25=\", \"55+12=\", \"16+29=\", \"89-56=\", \"91-3=\", \"12+45=\", \"81+2=\", \"32+23=\", \"7+66=\", \"88-59=\", \"58-46=\", \"25-18=\", \"48+28=\", \"41+9=\", \"2+23=\", \"24-11=\", \"25+50=\", \"24+60=\", \"17+46=\", \"40+53=\", \"8+14=\", \"89-77=\", \"82+3=\", \"93-81=\", \"93-82=\", \"90-34=\", \"90-38=\", \"4+40=\", \"21+74=\", \"46+34=\", \"92-53=\", \"69-21=\", \"89-49=\", \"62+10=\", \"28+2=\", \"38+31=\", \"75+15=\", \"4+95=\", \"10+61=\", \"37-1=\", \"78-6=\", \"42-1=\", \"27+48=\", \"18+60=\", \"63-42=\", \"90-56=\", \"85-21=\", \"87+1=\", \"23+35=\", \"31+11=\", \"95-61=\", \"81-50=\", \"22+26=\", \"70-53=\", \"55+16=\", \"90-88=\", \"95-94=\", \"11+69=\", \"34+4=\", \"64-59=\", \"48+47=\", \"7+4=\", \"55+9=\", \"21-6=\", \"70+14=\", \"92-12=\", \"27+33=\", \"54-45=\", \"22+58=\", \"28+19=\", \"95-87=\", \"48-31=\", \"58-36=\", \"34+55=\", \"61+4=\", \"55-5=\", \"91-10=\", \"97-23=\", \"1+83=\", \"73-64=\", \"93-70=\", \"42-1=\", \"97-76=\", \"86-49=\", \"23+15=\", \"17+42=\", \"37-25=\", \"87+11=\", \"95-10=\", \"7+68=\", \"37+42=\", \"54+0=\", \"87-9=\", \"38+54=\", \"83-70=\", \"79+2=\", \"77-72=\", \"57-14=\"];\n\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\nconst columnCount = 5;\nlet valueIndex = 0;\nfor (let r = 0; r < rows.items.length; r++) {\n  for (let c = 0; c < columnCount; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = newValues[valueIndex];\n    valueIndex++;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line and every arithmetic-problem cell in the practice\n# table to the new values from the authoritative diff. Cell values are\n# applied positionally (row-major, 5 columns per row) since some old\n# expressions repeat (e.g. \"63+4=\", \"9+84=\") and a naive Find/Replace could\n# not tell the duplicates apart.\n\n$d = $word.ActiveDocument\n\n# --- 1. Update the date paragraph ------------------------------------\n$d.Paragraphs.Item(1).Range.Text = \"2025-05-25 Sunday\"\n\n# --- 2. Update every cell in the practice table -----------------------\n$newValues = @(\"84-31=\", \"19+33=\", \"85-25=\", \"55+12=\", \"16+29=\", \"89-56=\", \"91-3=\", \"12+45=\", \"81+2=\", \"32+23=\", \"7+66=\", \"88-59=\", \"58-46=\", \"25-18=\", \"48+28=\", \"41+9=\", \"2+23=\", \"24-11=\", \"25+50=\", \"24+60=\", \"17+46=\", \"40+53=\", \"8+14=\", \"89-77=\", \"82+3=\", \"93-81=\", \"93-82=\", \"90-34=\", \"90-38=\", \"4+40=\", \"21+74=\", \"46+34=\", \"92-53=\", \"69-21=\", \"89-49=\", \"62+10=\", \"28+2=\", \"38+31=\", \"75+15=\", \"4+95=\", \"10+61=\", \"37-1=\", \"78-6=\", \"42-1=\", \"27+48=\", \"18+60=\", \"63-42=\", \"90-56=\", \"85-21=\", \"87+1=\", \"23+35=\", \"31+11=\", \"95-61=\", \"81-50=\", \"22+26=\", \"70-53=\", \"55+16=\", \"90-88=\", \"95-94=\", \"11+69=\", \"34+4=\", \"64-59=\", \"48+47=\", \"7+4=\", \"55+9=\", \"21-6=\", \"70+14=\", \"92-12=\", \"27+33=\", \"54-45=\", \"22+58=\", \"28+19=\", \"95-87=\", \"48-31=\", \"58-36=\", \"34+55=\", \"61+4=\", \"55-5=\", \"91-10=\", \"97-23=\", \"1+83=\", \"73-64=\", \"93-70=\", \"42-1=\", \"97-76=\", \"86-49=\", \"23+15=\", \"17+42=\", \"37-25=\", \"87+11=\", \"95-10=\", \"7+68=\", \"37+42=\", \"54+0=\", \"87-9=\", \"38+54=\", \"83-70=\", \"79+2=\", \"77-72=\", \"57-14=\")\n\n$table = $d.Tables.Item(1)\n$rowCount = $table.Rows.Count\n$columnCount = $table.Columns.Count\n\n$valueIndex = 0\nfor ($r = 1; $r -le $rowCount; $r++) {\n  for ($c = 1; $c -le $columnCount; $c++) {\n    $cell = $table.Cell($r, $c)\n    $cell.Range.Text = $newValues[$valueIndex]\n    $valueIndex = $valueIndex + 1\n  }\n}\n"}
